$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (shared strings are renumbered as a side-effect
# of Excel rewriting the sharedStrings table when these values change)
$ws.Range("A1").Value = "Row ID"
$ws.Range("C1").Value = "Task"
$ws.Range("E1").Value = "Start Date"
$ws.Range("F1").Value = "End Date"

# Update the active selection in the sheet view
$ws.Activate()
$ws.Range("F2").Select()
